$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump Version (B3) and concept Count (B22) ---
# Force text number-format first so "1.8.11" / "14" are not auto-parsed as a date/number
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").NumberFormat = "@"
$meta.Range("B3").Value = "1.8.11"
$meta.Range("B22").NumberFormat = "@"
$meta.Range("B22").Value = "14"

# --- Concepts sheet: rewrite concept table, add RNPI (13) and OTRO example rows ---
$concepts = $wb.Worksheets.Item("Concepts")

# Extend the formatted/bordered table from row 13 down to the two new rows (14-15)
# so they pick up the same border/alignment style as the rest of the table.
$concepts.Range("A13:D13").Copy()
$concepts.Range("A14:D15").PasteSpecial(-4122)
$concepts.Application.CutCopyMode = $false

# Columns A (Level) and B (Code) hold values that look numeric ("1", "01" .. "14");
# force them to Text so they keep leading zeros / do not turn into numbers.
$concepts.Range("A2:B15").NumberFormat = "@"

$concepts.Range("A2").Value = "1"
$concepts.Range("B2").Value = "01"
$concepts.Range("C2").Value = "RUN"
$concepts.Range("D2").Value = "Rol Único Nacional"
$concepts.Range("A3").Value = "1"
$concepts.Range("B3").Value = "02"
$concepts.Range("C3").Value = "RUN Provisorio"
$concepts.Range("D3").Value = "RUN provisorio (Artículo 44)"
$concepts.Range("A4").Value = "1"
$concepts.Range("B4").Value = "03"
$concepts.Range("C4").Value = "RUN Madre"
$concepts.Range("D4").Value = "RUN Madre (para recién nacido)"
$concepts.Range("A5").Value = "1"
$concepts.Range("B5").Value = "04"
$concepts.Range("C5").Value = "Número Folio"
$concepts.Range("D5").Value = "Número Folio Comprobante de Parto chileno"
$concepts.Range("A6").Value = "1"
$concepts.Range("B6").Value = "05"
$concepts.Range("C6").Value = "PPN"
$concepts.Range("D6").Value = "Pasaporte"
$concepts.Range("A7").Value = "1"
$concepts.Range("B7").Value = "06"
$concepts.Range("C7").Value = "Documento de identificación del país de origen"
$concepts.Range("D7").Value = "Documento de identificación del país de origen"
$concepts.Range("A8").Value = "1"
$concepts.Range("B8").Value = "07"
$concepts.Range("C8").Value = "Acta de nacimiento del país de origen"
$concepts.Range("D8").Value = "Acta de nacimiento del país de origen"
$concepts.Range("A9").Value = "1"
$concepts.Range("B9").Value = "08"
$concepts.Range("C9").Value = "NIP"
$concepts.Range("D9").Value = "Número de Identificación Provisorio (NIP)"
$concepts.Range("A10").Value = "1"
$concepts.Range("B10").Value = "09"
$concepts.Range("C10").Value = "NIC"
$concepts.Range("D10").Value = "Número Identificatorio para cotizar (NIC)"
$concepts.Range("A11").Value = "1"
$concepts.Range("B11").Value = "10"
$concepts.Range("C11").Value = "IPA"
$concepts.Range("D11").Value = "Identificación Provisoria del Apoderado (IPA)"
$concepts.Range("A12").Value = "1"
$concepts.Range("B12").Value = "11"
$concepts.Range("C12").Value = "IPE"
$concepts.Range("D12").Value = "Identificación Provisoria del Escolar (IPE)"
$concepts.Range("A13").Value = "1"
$concepts.Range("B13").Value = "12"
$concepts.Range("C13").Value = "Número de Ficha Clínica Sistema Local"
$concepts.Range("D13").Value = "Número de Ficha Clínica Sistema Local"
$concepts.Range("A14").Value = "1"
$concepts.Range("B14").Value = "13"
$concepts.Range("C14").Value = "RNPI"
$concepts.Range("D14").Value = "Registro Nacional de Prestadores Individuales"
$concepts.Range("A15").Value = "1"
$concepts.Range("B15").Value = "14"
$concepts.Range("C15").Value = "OTRO"
$concepts.Range("D15").Value = "Otro tipo de identificador"
